$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: "time_taken", styled like the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell F2: the recorded time_taken value for this panel
$ws.Range("F2").Value = "2021-10-05 13:42:08.434045"
